$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" values must stay as text cells, like the rest of the
# column (the sheet stores them as plain inline strings, not numbers).
# Briefly mark the cell as Text so the literal isn't coerced into a float,
# then restore the original "Normal" style so no stray formatting is left
# behind once the text value has been committed.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D18","D19","D20","D21","D23","D26","D40","D41","D42","D43","D44","D45","D46","D48","D49","D50")

$priceValues = @{
    "D2"  = "281.15"
    "D3"  = "20.67"
    "D4"  = "6.208"
    "D5"  = "0.06169"
    "D6"  = "3.583"
    "D7"  = "1.515"
    "D8"  = "6.567"
    "D9"  = "0.8186"
    "D10" = "0.01385"
    "D11" = "0.1641"
    "D12" = "0.08389"
    "D13" = "0.03538"
    "D14" = "0.03214"
    "D15" = "0.09136"
    "D16" = "3.716"
    "D18" = "0.04708"
    "D19" = "0.006416"
    "D20" = "0.006170"
    "D21" = "0.001069"
    "D23" = "3.779"
    "D26" = "0.1250"
    "D40" = "0.04705"
    "D41" = "0.007196"
    "D42" = "0.1101"
    "D43" = "0.003400"
    "D44" = "0.01102"
    "D45" = "0.00006594"
    "D46" = "0.00000000750"
    "D48" = "0.002857"
    "D49" = "0.00001901"
    "D50" = "0.01241"
}

foreach ($ref in $priceCells) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceValues[$ref]
    $cell.Style = "Normal"
}

# Plain text columns (Coin name / Link / Volume label) - swap rows 7/8 and 42/43,
# and refresh the Price-derived label text in column E for the other changed rows.
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E7").Value = "6FTXTokenFTT"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("E8").Value = "7KuCoinTokenKCS"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
